# Add new localization rows for weather-related strings (humidity, wind,
# and the weatherXxx / weatherXxxDesc key+value pairs) below the existing
# welcome/title rows, in the same order the author typed them so the
# shared-string table comes out in the same sequence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# humidity / wind
$ws.Range("A4").Value = "humidity"
$ws.Range("B4").Value = "Humidity"
$ws.Range("A5").Value = "wind"
$ws.Range("B5").Value = "Wind"

# Sunny
$ws.Range("A6").Value = "weatherSunny"
$ws.Range("B6").Value = "Sunny"
$ws.Range("A7").Value = "weatherSunnyDesc"
$ws.Range("B7").Value = "It's always sunny."

# Mostly Cloudy (keys typed first, then values)
$ws.Range("A10").Value = "weatherMostlyCloudy"
$ws.Range("A11").Value = "weatherMostlyCloudyDesc"
$ws.Range("B10").Value = "Mostly Cloudy"
$ws.Range("B11").Value = "Clouds, clouds everywhere."

# Partly Sunny (keys typed first, then values)
$ws.Range("A8").Value = "weatherPartlySunny"
$ws.Range("A9").Value = "weatherPartlySunnyDesc"
$ws.Range("B8").Value = "Partly Sunny"
$ws.Range("B9").Value = "Sunny'ish."

# Clear
$ws.Range("A14").Value = "weatherClear"
$ws.Range("B14").Value = "Clear"
$ws.Range("A15").Value = "weatherClearDesc"
$ws.Range("B15").Value = "No clouds allowed."

# Cloudy (keys typed first, then values)
$ws.Range("A12").Value = "weatherCloudy"
$ws.Range("A13").Value = "weatherCloudyDesc"
$ws.Range("B12").Value = "Cloudy"
$ws.Range("B13").Value = "Just clouds."

# Column widths were widened to fit the new, longer key/value strings.
$ws.Columns.Item(1).ColumnWidth = 34.6
$ws.Columns.Item(2).ColumnWidth = 84.6
$ws.Columns.Item(3).ColumnWidth = 19.45

# Selection moved to B12 as the last-edited cell.
$ws.Range("B12").Select() | Out-Null
